# "fix the problem of excel location"
#
# The SignIn sheet's Url cell (A2) pointed at a dead/obsolete address
# (http://www.skillswap.pro/Home). Point it at the correct test-environment
# location (http://192.168.99.100:5000/) and keep the hyperlink + displayed
# text + cell style in sync with that change.

$wb = $excel.ActiveWorkbook

$newUrl = "http://192.168.99.100:5000/"

$wsSignIn = $wb.Worksheets.Item("SignIn")
$urlCell = $wsSignIn.Range("A2")

# Drop the stale hyperlink (and its relationship) before re-pointing it,
# otherwise Excel keeps the old target around as an orphaned relationship.
$urlCell.Hyperlinks.Delete()

# The cell shows the URL as literal text, so update that too.
$urlCell.Value = $newUrl

# Re-create the hyperlink against the fixed location.
$wsSignIn.Hyperlinks.Add($urlCell, $newUrl)

# Hyperlinks.Add() reset the cell to a raw "applied" hyperlink font; put it
# back on the workbook's shared "Hyperlink" cell style like it originally was.
$urlCell.Style = "Hyperlink"

# --- incidental UI state saved alongside the fix ---
$wsSignUp = $wb.Worksheets.Item("SignUp")
[void]$wsSignUp.Range("B32").Select()

[void]$wsSignIn.Range("A2").Select()

# Restore ShareSkill as the active tab (selecting cells above switches sheets).
$wsShareSkill = $wb.Worksheets.Item("ShareSkill")
[void]$wsShareSkill.Activate()
